$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.172.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.647.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.541"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.138.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.274.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.545"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.785.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0974"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "
